# Update TPM-derived values in row 2 to reflect the new scaling factor.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.03635433333333333
$ws.Range("N2").Value = 0.109063
$ws.Range("Q2").Value = 0.04651513925588888
$ws.Range("R2").Value = 0.418636253303
